# XForm Demo: shift the slide-6 flow diagram to the right (every shape
# except the "Title 1" placeholder moves +746620 EMU on the X axis) and
# rename the "Official Dataset [Cloud Storage]" label to "Production
# Dataset [Cloud Storage]" on both slide 6 and slide 7.

# Shape.Left/.Top round-trip through a 32-bit float and the EMU<-points
# back-conversion truncates rather than rounds, so asking for the exact
# target EMU value can land 1 EMU low. Nudging by half an EMU before the
# unit conversion keeps the truncation from biting while staying far
# below the smallest meaningful distance (1 pt = 12700 EMU).
function EmuToPt($emu) {
    return ($emu + 0.5) / 914400 * 72
}

$p = $ppt.ActivePresentation

$slide6 = $p.Slides.Item(6)

$newLeftEmu = @{
    "Group 94" = 2738060
    "Group 86" = 2470946
    "Group 18" = 2225624
    "Group 19" = 2225624
    "Straight Arrow Connector 22" = 2926105
    "Straight Arrow Connector 77" = 3947519
    "Straight Arrow Connector 80" = 4438165
    "Straight Arrow Connector 84" = 5217716
    "Straight Arrow Connector 85" = 5708362
    "Right Arrow 115" = 999974
    "Rounded Rectangle 116" = 6939614
    "Rounded Rectangle 117" = 6939612
    "Rounded Rectangle 118" = 6939612
    "Straight Arrow Connector 119" = 5975503
}

for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    $target = $newLeftEmu[$shp.Name]
    if ($target -ne $null) {
        $shp.Left = EmuToPt($target)
    }
}

# --- Text updates: "Official Dataset [Cloud Storage]" -> "Production
#     Dataset [Cloud Storage]" inside "Group 18" / "Rectangle 16" on both
#     slide 6 and slide 7.
$grp6 = $slide6.Shapes.Item("Group 18")
$grp6.GroupItems.Item("Rectangle 16").TextFrame.TextRange.Text = "Production Dataset [Cloud Storage]"

$slide7 = $p.Slides.Item(7)
$grp7 = $slide7.Shapes.Item("Group 18")
$grp7.GroupItems.Item("Rectangle 16").TextFrame.TextRange.Text = "Production Dataset [Cloud Storage]"
